# Bates Hotel workbook - datos.xlsx edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# D3 (9) and D4 (55) are deleted outright (cell contents cleared, not just set to 0)
$ws.Range("D3").ClearContents() | Out-Null
$ws.Range("D4").ClearContents() | Out-Null

# E2 picks up an explicit (re-applied) font formatting, which creates and
# references a new cell style for that cell
$ws.Range("E2").Font.Name = "Calibri"

# The sheet's active/selected cell moves from E4 to D4
$ws.Range("D4").Select() | Out-Null
